# Apply the Sprint 2 documentation updates:
#   Sprint No.: 1 -> 2
#   Review Date: 02/09/18 -> 02/21/18

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 2, Cell 4 holds the Sprint No. value ("1").
$sprintCell = $tbl.Cell(2, 4)
$sprintCell.Range.Find.Execute("1", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "2", 1)

# Row 3, Cell 2 holds the Review Date value ("02/09/18").
$dateCell = $tbl.Cell(3, 2)
$dateCell.Range.Find.Execute("02/09/18", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "02/21/18", 1)
